$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3875
$ws.Range("J17").Value = 3875
$ws.Range("L17").Value = 11625
$ws.Range("N17").Value = -11961
$ws.Range("H18").Value = 1820
$ws.Range("I18").Value = 1820
$ws.Range("K18").Value = 1820
$ws.Range("M18").Value = -1536
$ws.Range("H52").Value = 331.75
$ws.Range("I52").Value = 34.75
$ws.Range("J52").Value = 628.75
$ws.Range("K52").Value = 104.25
$ws.Range("L52").Value = 1886.25
$ws.Range("M52").Value = 55.75
$ws.Range("N52").Value = -2206.25
$ws.Range("H53").Value = 341.92307
$ws.Range("I53").Value = 489.8
$ws.Range("J53").Value = 249.5
$ws.Range("K53").Value = 489.8
$ws.Range("L53").Value = 249.5
$ws.Range("M53").Value = 147.2
$ws.Range("N53").Value = -1523.5
$ws.Range("H64").Value = 3300
$ws.Range("I64").Value = 3600
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 3600
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -3352
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3300
$ws.Range("I67").Value = 3600
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 3600
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -2742
$ws.Range("N67").Value = -4916
$ws.Range("H75").Value = 53041.4
$ws.Range("J75").Value = 53041.4
$ws.Range("L75").Value = 53041.4
$ws.Range("N75").Value = -54913.4
$ws.Range("H78").Value = 53041.4
$ws.Range("J78").Value = 53041.4
$ws.Range("L78").Value = 159124.2
$ws.Range("N78").Value = -168484.2
$ws.Range("H86").Value = 3289
$ws.Range("I86").Value = 1859.4
$ws.Range("J86").Value = 4083.2222
$ws.Range("K86").Value = 1859.4
$ws.Range("L86").Value = 4083.2222
$ws.Range("M86").Value = -736.4000000000001
$ws.Range("N86").Value = -6329.2222
$ws.Range("H89").Value = 3289
$ws.Range("I89").Value = 1859.4
$ws.Range("J89").Value = 4083.2222
$ws.Range("K89").Value = 9297
$ws.Range("L89").Value = 20416.111
$ws.Range("M89").Value = -3681
$ws.Range("N89").Value = -31648.111
$ws.Range("H116").Value = 5660
$ws.Range("I116").Value = 3990
$ws.Range("K116").Value = 3990
$ws.Range("M116").Value = -548
$ws.Range("H135").Value = 1288.3636
$ws.Range("I135").Value = 1019.2222
$ws.Range("K135").Value = 9172.9998
$ws.Range("M135").Value = -6637.9998
$ws.Range("H137").Value = 9710.888999999999
$ws.Range("I137").Value = 7799.75
$ws.Range("K137").Value = 23399.25
$ws.Range("M137").Value = -20849.25
$ws.Range("H138").Value = 7739.108
$ws.Range("I138").Value = 7545.2173
$ws.Range("J138").Value = 7826.549
$ws.Range("K138").Value = 22635.6519
$ws.Range("L138").Value = 23479.647
$ws.Range("M138").Value = -17495.6519
$ws.Range("N138").Value = -33759.647

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20305.945
$ws.Range("J32").Value = 28442.445
$ws.Range("L32").Value = 28442.445
$ws.Range("N32").Value = -29016.445
$ws.Range("H74").Value = 2461.75
$ws.Range("I74").Value = 978.7692
$ws.Range("K74").Value = 978.7692
$ws.Range("M74").Value = -104.7692
$ws.Range("H77").Value = 2461.75
$ws.Range("I77").Value = 978.7692
$ws.Range("K77").Value = 4893.846
$ws.Range("M77").Value = -525.8459999999995
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 10000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = ""
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -10812
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 10000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = ""
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -12808
$ws.Range("H132").Value = 2847.2
$ws.Range("I132").Value = 1744.7693
$ws.Range("K132").Value = 5234.3079
$ws.Range("M132").Value = -2704.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2837.7222
$ws.Range("I20").Value = 2672.111
$ws.Range("J20").Value = 3003.3333
$ws.Range("K20").Value = 2672.111
$ws.Range("L20").Value = 3003.3333
$ws.Range("M20").Value = -2425.111
$ws.Range("N20").Value = -3497.3333
$ws.Range("H68").Value = 55000
$ws.Range("J68").Value = 55000
$ws.Range("L68").Value = 55000
$ws.Range("N68").Value = -56622
$ws.Range("H71").Value = 55000
$ws.Range("J71").Value = 55000
$ws.Range("L71").Value = 165000
$ws.Range("N71").Value = -173112
$ws.Range("H80").Value = 355
$ws.Range("J80").Value = 403.4
$ws.Range("L80").Value = 403.4
$ws.Range("N80").Value = -2399.4
$ws.Range("H83").Value = 355
$ws.Range("J83").Value = 403.4
$ws.Range("L83").Value = 2017
$ws.Range("N83").Value = -12001
$ws.Range("H86").Value = 5454.3335
$ws.Range("I86").Value = 4778
$ws.Range("K86").Value = 4778
$ws.Range("M86").Value = -3655
$ws.Range("H89").Value = 5454.3335
$ws.Range("I89").Value = 4778
$ws.Range("K89").Value = 23890
$ws.Range("M89").Value = -18274
$ws.Range("H134").Value = 2988.6667
$ws.Range("I134").Value = 1496.1428
$ws.Range("K134").Value = 4488.428400000001
$ws.Range("M134").Value = -1953.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6654.9165
$ws.Range("I86").Value = 3170.3125
$ws.Range("J86").Value = 13624.125
$ws.Range("K86").Value = 3170.3125
$ws.Range("L86").Value = 13624.125
$ws.Range("M86").Value = -2047.3125
$ws.Range("N86").Value = -15870.125
$ws.Range("H89").Value = 6654.9165
$ws.Range("I89").Value = 3170.3125
$ws.Range("J89").Value = 13624.125
$ws.Range("K89").Value = 15851.5625
$ws.Range("L89").Value = 68120.625
$ws.Range("M89").Value = -10235.5625
$ws.Range("N89").Value = -79352.625
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = ""
$ws.Range("N130").Value = ""
$ws.Range("H132").Value = 2203.9412
$ws.Range("I132").Value = 2042.2858
$ws.Range("J132").Value = 2958.3333
$ws.Range("K132").Value = 6126.857400000001
$ws.Range("L132").Value = 8874.999899999999
$ws.Range("M132").Value = -3596.857400000001
$ws.Range("N132").Value = -13934.9999
$ws.Range("H138").Value = 70780
$ws.Range("J138").Value = 70780
$ws.Range("L138").Value = 70780
$ws.Range("N138").Value = -81060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 1248
$ws.Range("I41").Value = 746.6667
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 2240.0001
$ws.Range("L41").Value = 6000
$ws.Range("M41").Value = -1902.0001
$ws.Range("N41").Value = -6676
$ws.Range("H96").Value = 6999
$ws.Range("J96").Value = 6999
$ws.Range("L96").Value = 20997
$ws.Range("N96").Value = -25115
$ws.Range("H117").Value = 1587.1
$ws.Range("J117").Value = 2248.5
$ws.Range("L117").Value = 6745.5
$ws.Range("N117").Value = -13629.5
$ws.Range("H131").Value = 1521.1052
$ws.Range("I131").Value = 1159.6
$ws.Range("J131").Value = 1575.8788
$ws.Range("K131").Value = 3478.8
$ws.Range("L131").Value = 4727.636399999999
$ws.Range("M131").Value = 1561.2
$ws.Range("N131").Value = -14807.6364
$ws.Range("H137").Value = 7975.857
$ws.Range("I137").Value = 9000
$ws.Range("K137").Value = 27000
$ws.Range("M137").Value = -21900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("H132").Value = 4600.6665
$ws.Range("I132").Value = 3375.7144
$ws.Range("K132").Value = 10127.1432
$ws.Range("M132").Value = -7597.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = ""
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = ""
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = ""
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = ""
$ws.Range("H46").Value = 2916.5833
$ws.Range("I46").Value = 2833.3333
$ws.Range("J46").Value = 2999.8333
$ws.Range("K46").Value = 2833.3333
$ws.Range("L46").Value = 2999.8333
$ws.Range("M46").Value = -2645.3333
$ws.Range("N46").Value = -3375.8333
$ws.Range("H122").Value = 7099.8184
$ws.Range("I122").Value = 5442.5713
$ws.Range("K122").Value = 16327.7139
$ws.Range("M122").Value = -13877.7139
$ws.Range("H132").Value = 4540.2666
$ws.Range("I132").Value = 2300.5715
$ws.Range("K132").Value = 6901.7145
$ws.Range("M132").Value = -4371.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1713
$ws.Range("I122").Value = 1713
$ws.Range("K122").Value = 5139
$ws.Range("M122").Value = -2689
$ws.Range("H132").Value = 1440.909
$ws.Range("I132").Value = 613.6667
$ws.Range("K132").Value = 1841.0001
$ws.Range("M132").Value = 688.9999
